$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.989999999999997
$ws.Range("D4").Value = -6.831599999999996
$ws.Range("C7").Value = -13.78109999999999
$ws.Range("B8").Value = 6.265699999999997
$ws.Range("B10").Value = 5.158
$ws.Range("E10").Value = 16.34079999999999
$ws.Range("D11").Value = -7.682200000000003
$ws.Range("B12").Value = 4.478599999999995
$ws.Range("E12").Value = 17.82610000000001
$ws.Range("E13").Value = 16.27800000000001
$ws.Range("C14").Value = -13.94590000000001
$ws.Range("D14").Value = -8.257000000000003
$ws.Range("E14").Value = 16.00050000000001
$ws.Range("C15").Value = -14.06149999999998
$ws.Range("B18").Value = 6.558499999999996
$ws.Range("C18").Value = -13.80939999999999
$ws.Range("D18").Value = -8.914299999999999
$ws.Range("D19").Value = -9.312899999999992
$ws.Range("C20").Value = -12.11570000000001
$ws.Range("D21").Value = -8.321199999999994
$ws.Range("B25").Value = 5.772599999999998
$ws.Range("D27").Value = -8.473299999999998
$ws.Range("C29").Value = -11.68130000000001
$ws.Range("E29").Value = 17.17800000000001
$ws.Range("C30").Value = -13.0321
$ws.Range("C31").Value = -12.7557
$ws.Range("D31").Value = -8.675800000000002
$ws.Range("E32").Value = 16.1313
$ws.Range("C35").Value = -12.0615
$ws.Range("E35").Value = 16.7403
$ws.Range("B37").Value = 8.691600000000001
$ws.Range("D38").Value = -8.321499999999993
$ws.Range("C40").Value = -13.3889
$ws.Range("D42").Value = -8.659999999999997
$ws.Range("E43").Value = 17.5521
$ws.Range("C44").Value = -13.45889999999999
$ws.Range("D44").Value = -7.984999999999999
$ws.Range("D47").Value = -7.606600000000001
$ws.Range("E48").Value = 17.55920000000002
$ws.Range("E49").Value = 15.7736
$ws.Range("C50").Value = -13.563
$ws.Range("E50").Value = 16.4061
$ws.Range("E51").Value = 17.3875
$ws.Range("C54").Value = -13.0423
$ws.Range("B55").Value = 6.144000000000001
$ws.Range("D56").Value = -8.783300000000001
$ws.Range("E56").Value = 16.4523
$ws.Range("D58").Value = -8.219599999999994
$ws.Range("E61").Value = 16.2666
$ws.Range("D65").Value = -7.827199999999999
$ws.Range("B68").Value = 5.5168
$ws.Range("C68").Value = -11.7098
$ws.Range("E69").Value = 17.49510000000003
$ws.Range("E71").Value = 16.6953
$ws.Range("D73").Value = -7.685899999999999
$ws.Range("C76").Value = -12.9362
$ws.Range("B77").Value = 8.777100000000004
$ws.Range("B78").Value = 8.973400000000005
$ws.Range("B79").Value = 8.226300000000005
$ws.Range("E79").Value = 18.64190000000002
$ws.Range("B80").Value = 9.081000000000003
$ws.Range("B81").Value = 5.973100000000001
$ws.Range("E81").Value = 16.36829999999998
$ws.Range("B82").Value = 5.525100000000002
$ws.Range("B84").Value = 6.1575
$ws.Range("C87").Value = -13.90639999999998
$ws.Range("C88").Value = -13.78429999999999
$ws.Range("D90").Value = -8.215699999999998
$ws.Range("C92").Value = -13.1383
$ws.Range("D92").Value = -7.941799999999996
$ws.Range("E92").Value = 16.3654
$ws.Range("D94").Value = -6.573299999999998
$ws.Range("D95").Value = -7.597799999999998
$ws.Range("C96").Value = -12.92580000000001
$ws.Range("C98").Value = -12.06219999999999
$ws.Range("B101").Value = 8.9893
$ws.Range("C101").Value = -12.7799
$ws.Range("D101").Value = -7.591899999999995
$ws.Range("B102").Value = 8.532900000000001
$ws.Range("C102").Value = -12.9155
